# Insert three new columns before the existing "Terms Typically Offered"
# column (D), which pushes it to column G.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1:F1").EntireColumn.Insert()

# New header labels for the inserted columns.
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# Fill the new columns with "NA" for every data row (2-32).
$ws.Range("D2:F32").Value = "NA"
